# Update countries & provincias Spain
#
# Applies the 23-Abr-2020 19:22 data refresh to the "Pais" sheet:
#  - bumps the "last updated" timestamp in A1
#  - refreshes the case/death/recovered counters for several countries
#  - Sudafrica overtakes Luxemburgo & Egipto, and Kazajistan overtakes
#    Hungria, in the case-count ranking, so those rows swap identities
#    while keeping their row position (the table stays sorted by total
#    cases, column B, descending)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Last updated stamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 19:22"

# --- Estados Unidos (row 4) ----------------------------------------------
$ws.Range("B4").Value = 866105
$ws.Range("C4").Value = 17388
$ws.Range("D4").Value = 84812
$ws.Range("E4").Value = 732484
$ws.Range("G4").Value = 1150
$ws.Range("H4").Value = 48809

# --- Turquia (row 10) -----------------------------------------------------
$ws.Range("B10").Value = 101790
$ws.Range("C10").Value = 3116
$ws.Range("D10").Value = 18491
$ws.Range("E10").Value = 80808
$ws.Range("F10").Value = 1816
$ws.Range("G10").Value = 115
$ws.Range("H10").Value = 2491

# --- Rumania (row 35) ------------------------------------------------------
$ws.Range("E35").Value = 7073
$ws.Range("G35").Value = 21
$ws.Range("H35").Value = 545

# --- Sudafrica overtakes Luxemburgo & Egipto (rows 53-55) ------------------
$ws.Range("A53").Value = "Sudafrica"
$ws.Range("B53").Value = 3953
$ws.Range("C53").Value = 318
$ws.Range("D53").Value = 1473
$ws.Range("E53").Value = 2405
$ws.Range("F53").Value = 36
$ws.Range("G53").Value = 10
$ws.Range("H53").Value = 75

$ws.Range("A54").Value = "Luxemburgo"
$ws.Range("B54").Value = 3665
$ws.Range("C54").Value = 11
$ws.Range("D54").Value = 728
$ws.Range("E54").Value = 2854
$ws.Range("F54").Value = 27
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 83

$ws.Range("A55").Value = "Egipto"
$ws.Range("B55").Value = 3659
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 935
$ws.Range("E55").Value = 2448
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 276

# --- Kazajistan overtakes Hungria (rows 63-64) ------------------------------
$ws.Range("A63").Value = "Kazajistan"
$ws.Range("B63").Value = 2289
$ws.Range("C63").Value = 154
$ws.Range("D63").Value = 560
$ws.Range("E63").Value = 1709
$ws.Range("F63").Value = 29
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 20

$ws.Range("A64").Value = "Hungria"
$ws.Range("B64").Value = 2284
$ws.Range("C64").Value = 116
$ws.Range("D64").Value = 390
$ws.Range("E64").Value = 1655
$ws.Range("F64").Value = 61
$ws.Range("G64").Value = 14
$ws.Range("H64").Value = 239

# --- Principado de Andorra (row 93) ----------------------------------------
$ws.Range("D93").Value = 333
$ws.Range("E93").Value = 353
